# Auto-generated edit script applying scheduled market-data refresh to Sheets/Phantom_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the rows that changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5099.3335
$ws.Range("I32").Value = 4500.5
$ws.Range("J32").Value = 5578.4
$ws.Range("K32").Value = 4500.5
$ws.Range("L32").Value = 5578.4
$ws.Range("M32").Value = -4174.5
$ws.Range("N32").Value = -6230.4
$ws.Range("H40").Value = 2722.2222
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2825
$ws.Range("H51").Value = 10292.6
$ws.Range("I51").Value = 9627.571
$ws.Range("J51").Value = 10874.5
$ws.Range("K51").Value = 9627.571
$ws.Range("L51").Value = 10874.5
$ws.Range("M51").Value = -9143.571
$ws.Range("N51").Value = -11842.5
$ws.Range("H112").Value = 3201.5
$ws.Range("J112").Value = 3491
$ws.Range("L112").Value = 10473
$ws.Range("N112").Value = -12689
$ws.Range("H121").Value = 2282.6667
$ws.Range("J121").Value = 2282.6667
$ws.Range("L121").Value = 6848.000100000001
$ws.Range("N121").Value = -10342.0001
$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820
$ws.Range("H125").Value = 798.04346
$ws.Range("J125").Value = 798.04346
$ws.Range("L125").Value = 7182.39114
$ws.Range("N125").Value = -12102.39114
$ws.Range("H138").Value = 1835.8
$ws.Range("I138").Value = 1369
$ws.Range("J138").Value = 3036.1428
$ws.Range("K138").Value = 4107
$ws.Range("L138").Value = 9108.428400000001
$ws.Range("M138").Value = 1033
$ws.Range("N138").Value = -19388.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2541.3076
$ws.Range("I45").Value = 2136.4167
$ws.Range("K45").Value = 2136.4167
$ws.Range("M45").Value = -1759.4167
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H110").Value = 20756.1
$ws.Range("I110").Value = 17589.166
$ws.Range("K110").Value = 17589.166
$ws.Range("M110").Value = -15544.166
$ws.Range("H132").Value = 7939.222
$ws.Range("I132").Value = 7350.476
$ws.Range("K132").Value = 22051.428
$ws.Range("M132").Value = -19521.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H94").Value = 4571.0835
$ws.Range("I94").Value = 4294.6
$ws.Range("K94").Value = 4294.6
$ws.Range("M94").Value = -3843.6
$ws.Range("H99").Value = 1229.7273
$ws.Range("I99").Value = 962.7
$ws.Range("K99").Value = 962.7
$ws.Range("M99").Value = 535.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 638.8570999999999
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 2851
$ws.Range("I31").Value = 4337.3335
$ws.Range("J31").Value = 2107.8333
$ws.Range("K31").Value = 4337.3335
$ws.Range("L31").Value = 2107.8333
$ws.Range("M31").Value = -4042.3335
$ws.Range("N31").Value = -2697.8333
$ws.Range("H34").Value = 2851
$ws.Range("I34").Value = 4337.3335
$ws.Range("J34").Value = 2107.8333
$ws.Range("K34").Value = 4337.3335
$ws.Range("L34").Value = 2107.8333
$ws.Range("M34").Value = -4135.3335
$ws.Range("N34").Value = -2511.8333
$ws.Range("H58").Value = 4599.6
$ws.Range("I58").Value = 1999
$ws.Range("J58").Value = 5249.75
$ws.Range("K58").Value = 1999
$ws.Range("L58").Value = 5249.75
$ws.Range("M58").Value = -1796
$ws.Range("N58").Value = -5655.75
$ws.Range("H92").Value = 24666.666
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 24666.666
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 24666.666
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -29658.666
$ws.Range("H99").Value = 2106
$ws.Range("I99").Value = 1712
$ws.Range("K99").Value = 1712
$ws.Range("M99").Value = -214
$ws.Range("H105").Value = 1279.25
$ws.Range("I105").Value = 851.5
$ws.Range("K105").Value = 851.5
$ws.Range("M105").Value = 895.5
$ws.Range("H113").Value = 638.8570999999999
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2106
$ws.Range("I126").Value = 1712
$ws.Range("K126").Value = 5136
$ws.Range("M126").Value = -2666
$ws.Range("H132").Value = 7183
$ws.Range("I132").Value = 6626.625
$ws.Range("K132").Value = 19879.875
$ws.Range("M132").Value = -17349.875
$ws.Range("H136").Value = 4599.6
$ws.Range("I136").Value = 1999
$ws.Range("J136").Value = 5249.75
$ws.Range("K136").Value = 5997
$ws.Range("L136").Value = 15749.25
$ws.Range("M136").Value = -3447
$ws.Range("N136").Value = -20849.25
$ws.Range("H141").Value = 699848.7
$ws.Range("J141").Value = 699848.7
$ws.Range("L141").Value = 699848.7
$ws.Range("N141").Value = -710208.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 245
$ws.Range("J22").Value = 245
$ws.Range("L22").Value = 735
$ws.Range("N22").Value = -1073
$ws.Range("H27").Value = 245
$ws.Range("J27").Value = 245
$ws.Range("L27").Value = 735
$ws.Range("N27").Value = -939
$ws.Range("H80").Value = 3998.3333
$ws.Range("I80").Value = 3499.25
$ws.Range("J80").Value = 4996.5
$ws.Range("K80").Value = 10497.75
$ws.Range("L80").Value = 14989.5
$ws.Range("M80").Value = -9561.75
$ws.Range("N80").Value = -16861.5
$ws.Range("H83").Value = 3998.3333
$ws.Range("I83").Value = 3499.25
$ws.Range("J83").Value = 4996.5
$ws.Range("K83").Value = 31493.25
$ws.Range("L83").Value = 44968.5
$ws.Range("M83").Value = -26813.25
$ws.Range("N83").Value = -54328.5
$ws.Range("H109").Value = 3325
$ws.Range("I109").Value = 3325
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 9975
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -8935
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 1073.2
$ws.Range("I113").Value = 1082.8445
$ws.Range("K113").Value = 3248.5335
$ws.Range("M113").Value = -1078.5335
$ws.Range("H122").Value = 10617.3
$ws.Range("I122").Value = 25463.75
$ws.Range("J122").Value = 719.6667
$ws.Range("K122").Value = 229173.75
$ws.Range("L122").Value = 6477.0003
$ws.Range("M122").Value = -226723.75
$ws.Range("N122").Value = -11377.0003
$ws.Range("H131").Value = 1456.0769
$ws.Range("I131").Value = 1007.25
$ws.Range("J131").Value = 1655.5555
$ws.Range("K131").Value = 3021.75
$ws.Range("L131").Value = 4966.666499999999
$ws.Range("M131").Value = 2018.25
$ws.Range("N131").Value = -15046.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8999
$ws.Range("I70").Value = 7799
$ws.Range("K70").Value = 7799
$ws.Range("M70").Value = -7529
$ws.Range("H73").Value = 8999
$ws.Range("I73").Value = 7799
$ws.Range("K73").Value = 7799
$ws.Range("M73").Value = -6863
$ws.Range("H80").Value = 5909.625
$ws.Range("I80").Value = 5937.9165
$ws.Range("J80").Value = 5824.75
$ws.Range("K80").Value = 5937.9165
$ws.Range("L80").Value = 5824.75
$ws.Range("M80").Value = -4939.9165
$ws.Range("N80").Value = -7820.75
$ws.Range("H83").Value = 5909.625
$ws.Range("I83").Value = 5937.9165
$ws.Range("J83").Value = 5824.75
$ws.Range("K83").Value = 29689.5825
$ws.Range("L83").Value = 29123.75
$ws.Range("M83").Value = -24697.5825
$ws.Range("N83").Value = -39107.75
$ws.Range("H107").Value = 1002.38464
$ws.Range("J107").Value = 3393.3333
$ws.Range("L107").Value = 3393.3333
$ws.Range("N107").Value = -7233.3333
$ws.Range("H113").Value = 2935.1667
$ws.Range("I113").Value = 2734.3333
$ws.Range("J113").Value = 3537.6667
$ws.Range("K113").Value = 2734.3333
$ws.Range("L113").Value = 3537.6667
$ws.Range("M113").Value = -564.3332999999998
$ws.Range("N113").Value = -7877.6667
$ws.Range("H132").Value = 4247
$ws.Range("I132").Value = 4247
$ws.Range("K132").Value = 12741
$ws.Range("M132").Value = -10211

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 749.8570999999999
$ws.Range("I22").Value = 669.75
$ws.Range("J22").Value = 856.6667
$ws.Range("K22").Value = 669.75
$ws.Range("L22").Value = 856.6667
$ws.Range("M22").Value = -374.75
$ws.Range("N22").Value = -1446.6667
$ws.Range("H27").Value = 749.8570999999999
$ws.Range("I27").Value = 669.75
$ws.Range("J27").Value = 856.6667
$ws.Range("K27").Value = 669.75
$ws.Range("L27").Value = 856.6667
$ws.Range("M27").Value = -562.75
$ws.Range("N27").Value = -1070.6667
$ws.Range("H82").Value = 2693.2727
$ws.Range("I82").Value = 1826.3334
$ws.Range("J82").Value = 3018.375
$ws.Range("K82").Value = 1826.3334
$ws.Range("L82").Value = 3018.375
$ws.Range("M82").Value = -1465.3334
$ws.Range("N82").Value = -3740.375
$ws.Range("H85").Value = 2693.2727
$ws.Range("I85").Value = 1826.3334
$ws.Range("J85").Value = 3018.375
$ws.Range("K85").Value = 1826.3334
$ws.Range("L85").Value = 3018.375
$ws.Range("M85").Value = -578.3334
$ws.Range("N85").Value = -5514.375
$ws.Range("H132").Value = 4422.2383
$ws.Range("I132").Value = 3174.375
$ws.Range("J132").Value = 5190.154
$ws.Range("K132").Value = 9523.125
$ws.Range("L132").Value = 15570.462
$ws.Range("M132").Value = -6993.125
$ws.Range("N132").Value = -20630.462
